# Update "想去人数" (want-to-go count) values in column F for five events,
# on both the "展览" sheet and the "全部类型" sheet (which mirrors the
# same events at slightly different row numbers).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 765   # was 764
$ws1.Range("F3").Value = 24    # was 23
$ws1.Range("F7").Value = 3627  # was 3626
$ws1.Range("F9").Value = 4249  # was 4246
$ws1.Range("F11").Value = 1058 # was 1057

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 765   # was 764
$ws4.Range("F3").Value = 24    # was 23
$ws4.Range("F8").Value = 3627  # was 3626
$ws4.Range("F10").Value = 4249 # was 4246
$ws4.Range("F12").Value = 1058 # was 1057
